$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 38-59 (NSE:KAUSHALYA .. NSE:SAFARI) are removed entirely; the used
# range shrinks from A1:F59 down to A1:F37.
$ws.Range("A38:F59").EntireRow.Delete()

# Refresh the watch-list columns for rows 2-37 with the new ticker symbols.
# Only cells whose content actually changes are touched, so cells that were
# already (and remain) blank are left completely alone.
$ws.Range("B2").Value = "NSE:AGARIND"
$ws.Range("C2").Value = "NSE:BATAINDIA"
$ws.Range("D2").Value = "NSE:ICICIGI"
$ws.Range("F2").Value = "NSE:HDFCBANK"
$ws.Range("C3").Value = "NSE:NEOGEN"
$ws.Range("F3").Value = "NSE:ICICIGI"
$ws.Range("B4").Value = "NSE:AMIORG"
$ws.Range("C4").Value = "NSE:PNBGILTS"
$ws.Range("B5").Value = "NSE:ASKAUTOLTD"
$ws.Range("C5").Value = "NSE:RELINFRA"
$ws.Range("B6").Value = "NSE:CCHHL"
$ws.Range("C6").Value = ""   # was "NSE:INVENTURE"
$ws.Range("B7").Value = "NSE:CELLO"
$ws.Range("C7").Value = ""   # was "NSE:NEOGEN"
$ws.Range("B8").Value = "NSE:COSMOFIRST"
$ws.Range("C8").Value = ""   # was "NSE:RAJRATAN"
$ws.Range("B9").Value = "NSE:DCM"
$ws.Range("B10").Value = "NSE:DCMNVL"
$ws.Range("B11").Value = "NSE:DEEPAKFERT"
$ws.Range("B12").Value = "NSE:DHARMAJ"
$ws.Range("B13").Value = "NSE:DIAMINESQ"
$ws.Range("B14").Value = "NSE:DSSL"
$ws.Range("B15").Value = "NSE:ENDURANCE"
$ws.Range("B16").Value = "NSE:ENIL"
$ws.Range("B17").Value = "NSE:ESTER"
$ws.Range("B18").Value = "NSE:FAIRCHEMOR"
$ws.Range("B19").Value = "NSE:GLOBUSSPR"
$ws.Range("B20").Value = "NSE:GNA"
$ws.Range("B21").Value = "NSE:GREAVESCOT"
$ws.Range("B22").Value = "NSE:GUFICBIO"
$ws.Range("B23").Value = "NSE:HDFCBANK"
$ws.Range("B24").Value = "NSE:HEMIPROP"
$ws.Range("B25").Value = "NSE:HINDOILEXP"
$ws.Range("B26").Value = "NSE:IZMO"
$ws.Range("B27").Value = "NSE:JINDALPOLY"
$ws.Range("B28").Value = "NSE:KBCGLOBAL"
$ws.Range("B29").Value = "NSE:KIOCL"
$ws.Range("B30").Value = "NSE:LUXIND"
$ws.Range("B31").Value = "NSE:MEDICAMEQ"
$ws.Range("B32").Value = "NSE:MIDHANI"
$ws.Range("B33").Value = "NSE:MOTOGENFIN"
$ws.Range("B34").Value = "NSE:PLASTIBLEN"
$ws.Range("B35").Value = "NSE:POLYPLEX"
$ws.Range("B36").Value = "NSE:PONNIERODE"
$ws.Range("B37").Value = "NSE:PTCIL"
